$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "Periodo Mora" column value from 2507 to 2508 for all affected rows.
$ws.Range("E16:E21").Value = "2508"

# Update the "Valor Mora" for the second worker (row 17) from 17,000,000 to 9,000,000.
$ws.Range("G17").Value = 9000000
